$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The patient name in A2 changes from "Mr. Ratan Singh" to "Mr. Indra Singh".
# (This also causes the shared-strings table entry for the old name to be
# dropped and a new entry for "Mr. Indra Singh" to be appended, which is
# reflected automatically by the engine on save.)
$ws.Range("A2").Value = "Mr. Indra Singh"

# The sheet's active selection moves to A2 (saved sheetView selection).
$null = $ws.Range("A2").Select()
